$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, $value)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "27.753.25"
$ws.Range("E2").Value = "  +2.70%  "
Set-TextValue $ws.Range("D3") "1.722.59"
$ws.Range("E3").Value = "  +2.79%  "
$ws.Range("E4").Value = "  -0.75%  "
Set-TextValue $ws.Range("D5") "217.51"
$ws.Range("E5").Value = "  +1.09%  "
$ws.Range("E6").Value = "  +0.80%  "
Set-TextValue $ws.Range("D7") "0.994"
$ws.Range("E7").Value = "  -0.73%  "
Set-TextValue $ws.Range("D8") "24.08"
$ws.Range("E8").Value = "  +12.52%  "
$ws.Range("E9").Value = "  +4.56%  "
Set-TextValue $ws.Range("D10") "0.0631"
$ws.Range("E10").Value = "  +1.29%  "
Set-TextValue $ws.Range("D11") "0.0897"
$ws.Range("E11").Value = "  +0.93%  "
Set-TextValue $ws.Range("D12") "1.958.22"
$ws.Range("E12").Value = "  +2.36%  "
Set-TextValue $ws.Range("D13") "1.703.74"
$ws.Range("E13").Value = "  +1.55%  "
$ws.Range("E14").Value = "  +3.12%  "
Set-TextValue $ws.Range("D15") "0.565"
$ws.Range("E15").Value = "  +5.90%  "
Set-TextValue $ws.Range("D16") "67.90"
$ws.Range("E16").Value = "  +2.41%  "
Set-TextValue $ws.Range("D17") "27.759.04"
$ws.Range("E17").Value = "  +2.69%  "
Set-TextValue $ws.Range("D18") "242.36"
$ws.Range("E18").Value = "  +2.72%  "
Set-TextValue $ws.Range("D19") "7.98"
$ws.Range("E19").Value = "  -2.82%  "
$ws.Range("E20").Value = "  +1.68%  "
Set-TextValue $ws.Range("D21") "0.994"
$ws.Range("E21").Value = "  -0.62%  "
$ws.Range("E22").Value = "  +3.56%  "
Set-TextValue $ws.Range("D23") "9.73"
$ws.Range("E23").Value = "  +5.36%  "
$ws.Range("E24").Value = "  +0.50%  "
Set-TextValue $ws.Range("D25") "148.68"
$ws.Range("E25").Value = "  +0.40%  "
$ws.Range("E26").Value = "  +4.06%  "
Set-TextValue $ws.Range("D27") "16.58"
$ws.Range("E28").Value = "  +1.23%  "
$ws.Range("E29").Value = "  -0.69%  "
$ws.Range("E30").Value = "  +1.53%  "
$ws.Range("E31").Value = "  +0.83%  "
Set-TextValue $ws.Range("D32") "3.44"
$ws.Range("E32").Value = "  +2.03%  "
Set-TextValue $ws.Range("D33") "1.553.20"
$ws.Range("E33").Value = "  +0.74%  "
Set-TextValue $ws.Range("D34") "3.31"
$ws.Range("E34").Value = "  +4.42%  "
$ws.Range("E35").Value = "  -1.65%  "
Set-TextValue $ws.Range("D36") "0.969"
$ws.Range("E36").Value = "  +6.05%  "
Set-TextValue $ws.Range("D37") "0.616"
$ws.Range("E37").Value = "  +4.43%  "
$ws.Range("E38").Value = "  +0.41%  "
$ws.Range("E39").Value = "  +0.13%  "
Set-TextValue $ws.Range("D40") "1.07"
$ws.Range("E40").Value = "  +2.19%  "
Set-TextValue $ws.Range("D41") "71.41"
$ws.Range("E41").Value = "  +5.31%  "
Set-TextValue $ws.Range("D42") "5.86"
$ws.Range("E42").Value = "  +5.82%  "
Set-TextValue $ws.Range("D43") "0.994"
$ws.Range("E43").Value = "  -0.70%  "
$ws.Range("E44").Value = "  +1.53%  "
Set-TextValue $ws.Range("D45") "1.864.42"
$ws.Range("E45").Value = "  +2.47%  "
$ws.Range("E46").Value = "  +1.54%  "
$ws.Range("E47").Value = "  +9.68%  "
Set-TextValue $ws.Range("D48") "91.82"
$ws.Range("E48").Value = "  +1.65%  "
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue $ws.Range("D49") "0.0₆0109"
$ws.Range("E49").Value = "  +1.04%  "
Set-TextValue $ws.Range("D50") "8.36"
$ws.Range("E50").Value = "  +4.78%  "
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue $ws.Range("D51") "0.106"
$ws.Range("E51").Value = "  +2.46%  "
